# Update labels "Drink", "Eat", "Smoke" to their uppercase forms
# "DRINK", "EAT", "SMOKE" (column A of Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$map = @{
    "Drink" = "DRINK"
    "Eat"   = "EAT"
    "Smoke" = "SMOKE"
}

$colA = $ws.Columns.Item(1)
foreach ($cell in $colA.Cells) {
    if ($cell.Row -gt $lastRow) {
        break
    }
    $v = $cell.Value2
    if ($v -ne $null -and $map.ContainsKey([string]$v)) {
        $cell.Value = $map[[string]$v]
    }
}
